$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.607.71'
$ws.Range('E2').Value = '  -2.43%  '
$ws.Range('D3').Value = '1.746.84'
$ws.Range('E3').Value = '  -3.43%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '322.26'
$ws.Range('E5').Value = '  -4.68%  '
$ws.Range('D6').Value = '0.9974'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4260'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -9.25%  '
$ws.Range('D8').Value = '0.3616'
$ws.Range('E8').Value = '  -6.45%  '
$ws.Range('D9').Value = '45.28'
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('D10').Value = '0.07473'
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range('E11').Value = '  -3.92%  '
$ws.Range('D12').Value = '0.9987'
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('D13').Value = '21.54'
$ws.Range('E13').Value = '  -4.55%  '
$ws.Range('D14').Value = '6.118'
$ws.Range('E14').Value = '  -3.97%  '
$ws.Range('D15').Value = '7.223'
$ws.Range('E15').Value = '  -3.28%  '
$ws.Range('D16').Value = '1.742.09'
$ws.Range('E16').Value = '  -3.77%  '
$ws.Range('D17').Value = '0.00001068'
$ws.Range('E17').Value = '  -2.69%  '
$ws.Range('D18').Value = '87.79'
$ws.Range('E18').Value = '  +7.32%  '
$ws.Range('D19').Value = '0.06228'
$ws.Range('E19').Value = '  -7.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9990'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '16.96'
$ws.Range('E21').Value = '  -3.41%  '
$ws.Range('E22').Value = '  -4.78%  '
$ws.Range('D23').Value = '0.5243'
$ws.Range('E23').Value = '  -6.41%  '
$ws.Range('D24').Value = '27.625.02'
$ws.Range('E24').Value = '  -2.36%  '
$ws.Range('E25').Value = '  -2.33%  '
$ws.Range('D26').Value = '2.317'
$ws.Range('E26').Value = '  -4.38%  '
$ws.Range('D27').Value = '20.48'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('D28').Value = '2.371'
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').Value = '151.71'
$ws.Range('E29').Value = '  -0.99%  '
$ws.Range('D30').Value = '1.941.14'
$ws.Range('E30').Value = '  -3.64%  '
$ws.Range('D31').Value = '1.217'
$ws.Range('E31').Value = '  -3.51%  '
$ws.Range('D32').Value = '126.99'
$ws.Range('E32').Value = '  -4.80%  '
$ws.Range('D33').Value = '5.711'
$ws.Range('E33').Value = '  -2.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09150'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.53%  '
$ws.Range('D35').Value = '3.673'
$ws.Range('E35').Value = '  -8.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.70'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.70%  '
$ws.Range('D37').Value = '0.02305'
$ws.Range('E37').Value = '  -2.28%  '
$ws.Range('D38').Value = '0.2138'
$ws.Range('E38').Value = '  -7.18%  '
$ws.Range('D39').Value = '5.096'
$ws.Range('E39').Value = '  -3.49%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.6449'
$ws.Range('E40').Value = '  -2.99%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '0.06094'
$ws.Range('E41').Value = '  -4.56%  '
$ws.Range('D42').Value = '1.191'
$ws.Range('E42').Value = '  -3.63%  '
$ws.Range('E43').Value = '  -5.09%  '
$ws.Range('D44').Value = '7.931'
$ws.Range('E44').Value = '  -5.61%  '
$ws.Range('D45').Value = '0.9973'
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.33%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.5909'
$ws.Range('E47').Value = '  -3.98%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').Value = '3.725'
$ws.Range('E48').Value = '  -3.40%  '
$ws.Range('D49').Value = '126.01'
$ws.Range('E49').Value = '  -4.01%  '
$ws.Range('D50').Value = '1.962'
$ws.Range('E50').Value = '  -3.98%  '
$ws.Range('D51').Value = '0.06878'
$ws.Range('E51').Value = '  -3.77%  '
